$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 11 and 12 represent two observation records whose contents were
# swapped between each other (only the columns that actually differ
# between the two rows need to be exchanged).
$cols = @("A","B","D","E","F","G","H","Q","R","AC","AX")

foreach ($col in $cols) {
    $addr11 = "$col`11"
    $addr12 = "$col`12"
    $v11 = $ws.Range($addr11).Value2
    $v12 = $ws.Range($addr12).Value2
    $ws.Range($addr11).Value2 = $v12
    $ws.Range($addr12).Value2 = $v11
}
